$wb = $excel.ActiveWorkbook

# --- Rename header labels on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after "Monthly Trend" (last sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "PO Forecast"

# --- Header row ---
$ws.Range("A1").Value = "ds"
$ws.Range("B1").Value = "PO_Forecast"
$ws.Range("C1").Value = "yhat_lower"
$ws.Range("D1").Value = "yhat_upper"

# --- Match formatting of the existing sheets: bold/centered header, date-formatted column A ---
$wsMonthly.Range("A1:B1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

$wsMonthly.Range("A2").Copy()
$ws.Range("A2:A62").PasteSpecial(-4122)

# --- Data rows ---
$data = @(
  @(45067.99999999999, 38, -21.56947898360227, 96.44638268498753),
  @(45074.99999999999, 39, -16.06459356404544, 95.18556641541537),
  @(45081.99999999999, 39, -20.19803211642097, 97.80807110026866),
  @(45095.99999999999, 40, -16.09108999546001, 100.6163304216431),
  @(45102.99999999999, 41, -19.84667411427937, 99.34368640042821),
  @(45109.99999999999, 41, -16.1365884355514, 95.31022937147294),
  @(45116.99999999999, 41, -17.94290576610734, 98.42043779368231),
  @(45130.99999999999, 42, -11.91684279798322, 102.172662283361),
  @(45137.99999999999, 43, -16.77364356507379, 103.9581034915831),
  @(45172.99999999999, 45, -11.93712681629538, 101.8745151253882),
  @(45186.99999999999, 46, -9.639574855981255, 102.2216114249332),
  @(45193.99999999999, 46, -15.70678260926194, 100.7383479527824),
  @(45207.99999999999, 47, -9.346439970004877, 106.1847718646503),
  @(45221.99999999999, 48, -7.307869304190483, 104.4167140253729),
  @(45228.99999999999, 48, -6.303075474523538, 103.1730568757035),
  @(45242.99999999999, 49, -9.236755500857511, 105.9315260144079),
  @(45249.99999999999, 50, -5.987107936795273, 106.2963876587962),
  @(45256.99999999999, 50, -11.46191746637973, 103.9168865864943),
  @(45263.99999999999, 51, -8.241464287666314, 109.9512557145808),
  @(45298.99999999999, 53, -4.224014481747234, 113.2421037069729),
  @(45305.99999999999, 53, -0.5559241870714017, 113.4613811018246),
  @(45319.99999999999, 54, -1.994716357237076, 112.0652877121648),
  @(45333.99999999999, 55, -2.451207204564707, 116.6600797456475),
  @(45347.99999999999, 56, 1.985975556923385, 115.7716936531899),
  @(45354.99999999999, 56, -6.049140606354721, 116.3945457049089),
  @(45361.99999999999, 57, -0.7246860502328782, 110.7092278456987),
  @(45382.99999999999, 58, 1.698959086035474, 119.9057435876909),
  @(45389.99999999999, 58, 4.071316238124396, 117.129643396759),
  @(45396.99999999999, 59, -1.048073926635135, 115.6060914759105),
  @(45403.99999999999, 59, 6.062160896900005, 120.9354392775677),
  @(45410.99999999999, 60, 2.874218304835999, 116.9392575602251),
  @(45424.99999999999, 60, 2.233159042940045, 115.0918554939968),
  @(45431.99999999999, 61, 6.548878522385951, 120.6526640017839),
  @(45445.99999999999, 62, 2.851959983385283, 118.2897138934654),
  @(45452.99999999999, 62, 2.069569326442187, 120.6554287670408),
  @(45459.99999999999, 63, 2.366917239463639, 120.854703511711),
  @(45466.99999999999, 63, 8.164207840526938, 120.9698386957511),
  @(45473.99999999999, 64, 6.194243017075834, 120.0893481820443),
  @(45480.99999999999, 64, 0.831088972430985, 121.3015605561995),
  @(45487.99999999999, 64, 8.331099515616705, 124.4604169716828),
  @(45501.99999999999, 65, 7.032114084276478, 122.0339016890793),
  @(45515.99999999999, 66, 5.568069167430245, 124.2054007853055),
  @(45522.99999999999, 67, 9.13164475247123, 124.5610635983838),
  @(45529.99999999999, 67, 8.064732614200201, 121.8236405104145),
  @(45536.99999999999, 67, 6.905198696228286, 124.7308239669705),
  @(45543.99999999999, 68, 6.581810012697339, 126.6836844656099),
  @(45550.99999999999, 68, 14.94358358270074, 127.9267940910617),
  @(45557.99999999999, 69, 8.25029294335641, 127.1160455275774),
  @(45564.99999999999, 69, 11.18692300559469, 127.2544530177933),
  @(45571.99999999999, 70, 15.60831058241384, 124.1868027571158),
  @(45585.99999999999, 70, 16.16044596345555, 128.4482216499758),
  @(45592.99999999999, 71, 11.40684178307739, 129.1380387993974),
  @(45599.99999999999, 71, 7.774063952032608, 126.4728553559534),
  @(45606.99999999999, 72, 11.22966488602484, 131.1600653146052),
  @(45613.99999999999, 72, 13.48131090476766, 130.6659494299773),
  @(45620.99999999999, 73, 16.23031821255336, 131.8403897092846),
  @(45627.99999999999, 73, 14.96563861913365, 129.4868940965985),
  @(45634.99999999999, 73, 21.7428572655664, 132.9733464224088),
  @(45641.99999999999, 74, 17.12252326700666, 136.3005435859707),
  @(45648.99999999999, 74, 14.96130071469526, 132.1847360694601),
  @(45655.99999999999, 75, 18.82254981178779, 133.5622926964765)
)
$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
